$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.21066472588013
$ws.Range("C2").Value = 8.530387320246627
$ws.Range("E2").Value = 20.59108797255688
$ws.Range("F2").Value = 39.09172823998971
$ws.Range("G2").Value = 29.12399235107401
$ws.Range("H2").Value = 13.98746694014338
$ws.Range("I2").Value = 19.50184430017446
$ws.Range("J2").Value = 7.75506064126036
$ws.Range("M2").Value = 19.01571478238296
$ws.Range("N2").Value = 17.13711488502166

$ws.Range("B3").Value = 11.66059121977759
$ws.Range("C3").Value = 7.99166715840295
$ws.Range("E3").Value = 20.57844340187529
$ws.Range("F3").Value = 38.97851533948202
$ws.Range("G3").Value = 28.90406832861767
$ws.Range("H3").Value = 14.01962747313627
$ws.Range("I3").Value = 19.58204174299138
$ws.Range("J3").Value = 7.777561774125193
$ws.Range("M3").Value = 18.84347314560693
$ws.Range("N3").Value = 17.20118367077483

$ws.Range("B4").Value = 11.3111879825029
$ws.Range("C4").Value = 7.640567978458895
$ws.Range("E4").Value = 20.57383100254083
$ws.Range("F4").Value = 38.91993182275736
$ws.Range("G4").Value = 28.78145858194685
$ws.Range("H4").Value = 14.04279225061137
$ws.Range("I4").Value = 19.63663790764872
$ws.Range("J4").Value = 7.792058413534594
$ws.Range("M4").Value = 18.7405308983133
$ws.Range("N4").Value = 17.24239196273313

$ws.Range("B5").Value = 11.16607899881047
$ws.Range("C5").Value = 7.492382370065542
$ws.Range("E5").Value = 20.57274621609064
$ws.Range("F5").Value = 38.89882139030694
$ws.Range("G5").Value = 28.7346673238357
$ws.Range("H5").Value = 14.05308864210909
$ws.Range("I5").Value = 19.66022640801038
$ws.Range("J5").Value = 7.798137645915915
$ws.Range("M5").Value = 18.69932873497297
$ws.Range("N5").Value = 17.25965629369056

$ws.Range("B6").Value = 11.14182582389331
$ws.Range("C6").Value = 7.467467727736671
$ws.Range("E6").Value = 20.57261415395496
$ws.Range("F6").Value = 38.89548324348378
$ws.Range("G6").Value = 28.72709060452253
$ws.Range("H6").Value = 14.05484999394166
$ws.Range("I6").Value = 19.6642240132883
$ws.Range("J6").Value = 7.799157485670078
$ws.Range("M6").Value = 18.69253343915911
$ws.Range("N6").Value = 17.26255154973789

$ws.Range("B7").Value = 11.30924172859712
$ws.Range("C7").Value = 7.638590181856324
$ws.Range("E7").Value = 20.57381315177275
$ws.Range("F7").Value = 38.91963591698652
$ws.Range("G7").Value = 28.7808146315914
$ws.Range("H7").Value = 14.04292764726633
$ws.Range("I7").Value = 19.63695061282379
$ws.Range("J7").Value = 7.792139704148174
$ws.Range("M7").Value = 18.73997215280965
$ws.Range("N7").Value = 17.2426228840818

$ws.Range("B8").Value = 12.02353385423565
$ws.Range("C8").Value = 8.348853792971669
$ws.Range("E8").Value = 20.58607556176139
$ws.Range("F8").Value = 39.0504336766458
$ws.Range("G8").Value = 29.04561557186657
$ws.Range("H8").Value = 13.99784475455236
$ws.Range("I8").Value = 19.52838124546256
$ws.Range("J8").Value = 7.762678039888056
$ws.Range("M8").Value = 18.95576595419515
$ws.Range("N8").Value = 17.15881859484578

$ws.Range("B9").Value = 13.32407251961459
$ws.Range("C9").Value = 9.580583026158397
$ws.Range("E9").Value = 20.63501116161179
$ws.Range("F9").Value = 39.39286050111075
$ws.Range("G9").Value = 29.66094947843019
$ws.Range("H9").Value = 13.936682253188
$ws.Range("I9").Value = 19.3582430147829
$ws.Range("J9").Value = 7.710281727755159
$ws.Range("M9").Value = 19.3994618106582
$ws.Range("N9").Value = 17.00924831225627

$ws.Range("B10").Value = 14.20990307078912
$ws.Range("C10").Value = 10.38773724227614
$ws.Range("E10").Value = 20.68596305629316
$ws.Range("F10").Value = 39.69555373978393
$ws.Range("G10").Value = 30.16773095172843
$ws.Range("H10").Value = 13.90851307908238
$ws.Range("I10").Value = 19.25968989596816
$ws.Range("J10").Value = 7.675030270434299
$ws.Range("M10").Value = 19.73546718626993
$ws.Range("N10").Value = 16.90826935687071

$ws.Range("B11").Value = 14.59637134260016
$ws.Range("C11").Value = 10.73382234722357
$ws.Range("E11").Value = 20.71235459101725
$ws.Range("F11").Value = 39.84401643385373
$ws.Range("G11").Value = 30.40920862286668
$ws.Range("H11").Value = 13.8993680333964
$ws.Range("I11").Value = 19.2206764602247
$ws.Range("J11").Value = 7.659690557671786
$ws.Range("M11").Value = 19.88995820374026
$ws.Range("N11").Value = 16.86424637181249

$ws.Range("B12").Value = 14.74025193281901
$ws.Range("C12").Value = 10.86185364935779
$ws.Range("E12").Value = 20.72280596064698
$ws.Range("F12").Value = 39.90175015473029
$ws.Range("G12").Value = 30.50213682105981
$ws.Range("H12").Value = 13.89643462755601
$ws.Range("I12").Value = 19.20674605146117
$ws.Range("J12").Value = 7.653981393161915
$ws.Range("M12").Value = 19.94864989441476
$ws.Range("N12").Value = 16.84784961761227

$ws.Range("B13").Value = 14.70937570707849
$ws.Range("C13").Value = 10.834414131134
$ws.Range("E13").Value = 20.72053480236037
$ws.Range("F13").Value = 39.88924939598203
$ws.Range("G13").Value = 30.48205841722339
$ws.Range("H13").Value = 13.89704280496128
$ws.Range("I13").Value = 19.20970861581007
$ws.Range("J13").Value = 7.655206539493327
$ws.Range("M13").Value = 19.93600191040917
$ws.Range("N13").Value = 16.85136879493138

$ws.Range("B14").Value = 14.60825840217258
$ws.Range("C14").Value = 10.74441609900577
$ws.Range("E14").Value = 20.71320529080581
$ws.Range("F14").Value = 39.84873608159169
$ws.Range("G14").Value = 30.41682461654827
$ws.Range("H14").Value = 13.89911607542659
$ws.Range("I14").Value = 19.21951346350947
$ws.Range("J14").Value = 7.659218866852872
$ws.Range("M14").Value = 19.89478327720704
$ws.Range("N14").Value = 16.86289192093204

$ws.Range("B15").Value = 14.54599739114952
$ws.Range("C15").Value = 10.6888962424235
$ws.Range("E15").Value = 20.7087751863778
$ws.Range("F15").Value = 39.82411663141377
$ws.Range("G15").Value = 30.37705787494481
$ws.Range("H15").Value = 13.90045504067418
$ws.Range("I15").Value = 19.22562920321002
$ws.Range("J15").Value = 7.661689496557213
$ws.Range("M15").Value = 19.86955894478398
$ws.Range("N15").Value = 16.86998578710965

$ws.Range("B16").Value = 14.18430605518867
$ws.Range("C16").Value = 10.36469611099208
$ws.Range("E16").Value = 20.68430258235398
$ws.Range("F16").Value = 39.68606515128992
$ws.Range("G16").Value = 30.15216258459393
$ws.Range("H16").Value = 13.9091847403222
$ws.Range("I16").Value = 19.26235715604489
$ws.Range("J16").Value = 7.67604672776335
$ws.Range("M16").Value = 19.72539986007092
$ws.Range("N16").Value = 16.91118474484379

$ws.Range("B17").Value = 13.95812222545018
$ws.Range("C17").Value = 10.16041381044151
$ws.Range("E17").Value = 20.67010927739535
$ws.Range("F17").Value = 39.60410919514023
$ws.Range("G17").Value = 30.01693776935132
$ws.Range("H17").Value = 13.91548141279504
$ws.Range("I17").Value = 19.28638351468798
$ws.Range("J17").Value = 7.6850324318014
$ws.Range("M17").Value = 19.63734992438202
$ws.Range("N17").Value = 16.93694798745776

$ws.Range("B18").Value = 13.82647870929384
$ws.Range("C18").Value = 10.0409299793639
$ws.Range("E18").Value = 20.66224841042378
$ws.Range("F18").Value = 39.55798587121732
$ws.Range("G18").Value = 29.94019528102434
$ws.Range("H18").Value = 13.91944835720898
$ws.Range("I18").Value = 19.30075016855831
$ws.Range("J18").Value = 7.690266343636165
$ws.Range("M18").Value = 19.58686313268712
$ws.Range("N18").Value = 16.95194644420033

$ws.Range("B19").Value = 13.78164362429949
$ws.Range("C19").Value = 10.00013342443622
$ws.Range("E19").Value = 20.65963898675836
$ws.Range("F19").Value = 39.54254472126271
$ws.Range("G19").Value = 29.91439199576761
$ws.Range("H19").Value = 13.92085073104381
$ws.Range("I19").Value = 19.30570827754951
$ws.Range("J19").Value = 7.692049733180749
$ws.Range("M19").Value = 19.56979759324841
$ws.Range("N19").Value = 16.95705563894768

$ws.Range("B20").Value = 13.98236087870636
$ws.Range("C20").Value = 10.18236549283213
$ws.Range("E20").Value = 20.67158887787182
$ws.Range("F20").Value = 39.61272867565398
$ws.Range("G20").Value = 30.03122615797671
$ws.Range("H20").Value = 13.91477537232814
$ws.Range("I20").Value = 19.28376917457085
$ws.Range("J20").Value = 7.684069105257683
$ws.Range("M20").Value = 19.64670705833272
$ws.Range("N20").Value = 16.93418681368833

$ws.Range("B21").Value = 14.63802661970051
$ws.Range("C21").Value = 10.77093268159735
$ws.Range("E21").Value = 20.71534576895667
$ws.Range("F21").Value = 39.86059501654938
$ws.Range("G21").Value = 30.43594575409625
$ws.Range("H21").Value = 13.89849271724709
$ws.Range("I21").Value = 19.21661060963981
$ws.Range("J21").Value = 7.658037648353895
$ws.Range("M21").Value = 19.90688542232003
$ws.Range("N21").Value = 16.85949987881868

$ws.Range("B22").Value = 15.05213189821266
$ws.Range("C22").Value = 11.13797648765519
$ws.Range("E22").Value = 20.74660750482243
$ws.Range("F22").Value = 40.03139769953165
$ws.Range("G22").Value = 30.70907110772034
$ws.Range("H22").Value = 13.89093884935883
$ws.Range("I22").Value = 19.17763586657484
$ws.Range("J22").Value = 7.641605212111813
$ws.Range("M22").Value = 20.07800890019353
$ws.Range("N22").Value = 16.81228295886792

$ws.Range("B23").Value = 14.83246180561841
$ws.Range("C23").Value = 10.94368702625946
$ws.Range("E23").Value = 20.72968038006026
$ws.Range("F23").Value = 39.93944315513389
$ws.Range("G23").Value = 30.56254021582267
$ws.Range("H23").Value = 13.89468737037344
$ws.Range("I23").Value = 19.19798542580867
$ws.Range("J23").Value = 7.650322543817227
$ws.Range("M23").Value = 19.98659344143348
$ws.Range("N23").Value = 16.83733796167933

$ws.Range("B24").Value = 13.97140758303713
$ws.Range("C24").Value = 10.17244748976969
$ws.Range("E24").Value = 20.67091901861958
$ws.Range("F24").Value = 39.60882870885379
$ws.Range("G24").Value = 30.02476325567601
$ws.Range("H24").Value = 13.9150934927304
$ws.Range("I24").Value = 19.28494939376904
$ws.Range("J24").Value = 7.684504413601188
$ws.Range("M24").Value = 19.64247627619285
$ws.Range("N24").Value = 16.9354345582575

$ws.Range("B25").Value = 12.9839426928405
$ws.Range("C25").Value = 9.264579247254993
$ws.Range("E25").Value = 20.61912469962488
$ws.Range("F25").Value = 39.29114386516326
$ws.Range("G25").Value = 29.48457581394294
$ws.Range("H25").Value = 13.95029438193905
$ws.Range("I25").Value = 19.39965220366695
$ws.Range("J25").Value = 7.723884146277543
$ws.Range("M25").Value = 19.2774871610418
$ws.Range("N25").Value = 17.04813959143857
